# Add a "2022-Q3" sheet (new quarter) right after "总计" and before "2022-Q2",
# with its fund-holdings table, and insert the matching summary row into
# the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: insert a new row 2 for 2022-Q3
#    and bump the index column (A) for every existing row by one.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows(2).Insert()
$summary.Range("B2:D2").ClearFormats()
$summary.Range("A3").Copy($summary.Range("A2"))

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 6
$summary.Cells.Item(2, 4).Value = 0.04

$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(6, 1).Value = 4

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet before the current "2022-Q2" tab
#    (tab #2) and clone the fund-table formatting from a sibling quarter
#    sheet so headers/index column keep the bold-centered-bordered style.
# ---------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item(2)
$formatSrc = $wb.Worksheets.Item(3)

$q3 = $wb.Worksheets.Add($refSheet)
$q3.Name = "2022-Q3"

$formatSrc.Range("A1:H7").Copy()
$q3.Range("A1").PasteSpecial(-4122)

# Header row
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Data rows: index, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$rows = @(
    @(0, "012675", "华融融泽6个月定开混合A", "1.27", "57.46", "1.30", "0.0165", 5),
    @(1, "515870", "嘉实中证先进制造100策略ETF", "0.36", "98.05", "2.47", "0.0089", 10),
    @(2, "080007", "长盛同鑫行业配置混合A", "0.20", "88.49", "2.32", "0.0046", 9),
    @(3, "012676", "华融融泽6个月定开混合C", "0.23", "57.46", "1.30", "0.0030", 5),
    @(4, "080015", "长盛中小盘精选混合", "0.13", "84.41", "2.23", "0.0029", 8),
    @(5, "010991", "长盛同鑫行业配置混合C", "0.02", "88.49", "2.32", "0.0005", 9)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $excelRow = $r + 2

    $q3.Cells.Item($excelRow, 1).Value = $row[0]

    # 基金代码 (fund code, col B) / 基金名称 (col C) / 基金规模 / 股票总仓位 /
    # 仓位占比 / 持有市值(亿元) (cols D-G) are stored as literal text (not
    # numbers) in the source workbook -- fund codes have significant leading
    # zeros (e.g. "012675") that numeric storage would drop. Force text
    # formatting before assignment, then drop back to the Normal style so no
    # stray number-format style lingers on the cell.
    foreach ($col in 2, 3, 4, 5, 6, 7) {
        $cell = $q3.Cells.Item($excelRow, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$col - 1]
        $cell.Style = "Normal"
    }

    $q3.Cells.Item($excelRow, 8).Value = $row[7]
}
